$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.962336003780365
$ws.Range("B1").Value = 1.898081064224243
$ws.Range("C1").Value = 4.893202304840088
$ws.Range("D1").Value = 2.230626583099365
$ws.Range("E1").Value = 0.5241938233375549
